# The "daftar_or_nah" column (C) values ("Yes" in every data row) are no
# longer needed - clear the data cells but keep the C1 header in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C50").ClearContents()

# Resize columns B and C to fit their (now unchanged/shorter) content.
$ws.Columns("B:C").AutoFit() | Out-Null

# Leave the active selection on D5, as in the saved file.
$ws.Range("D5").Select() | Out-Null
